# Added Presentation Project Folder
# Week 7 attendance data entry + tab/selection change to Week 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 7")

# Tuesday (row 3): hours per person, plus a remark in J3
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 4
$ws.Range("J3").Value = "Rilana is naar de tandarts"

# Wednesday (row 4)
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 2

# Thursday (row 5)
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4

# Friday (row 6)
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 8

# Geoorloofd (row 12) - Wednesday absence hours
$ws.Range("G12").Value = 4
$ws.Range("G12").HorizontalAlignment = -4108

# Make "Week 7" the active/selected sheet & tab, with J14 selected
$ws.Activate()
$ws.Range("J14").Select()
